# Updates cryptocurrency Price (D) and Volume(1h) (E) columns
# to match the latest scraped values (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.474.92'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').Value = '2.055.78'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '242.25'
$ws.Range('E5').Value = '  -1.84%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.662'
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '54.46'
$ws.Range('E8').Value = '  -5.29%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '58.57'
$ws.Range('E9').Value = '  -2.36%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.360'
$ws.Range('E10').Value = '  -5.63%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0750'
$ws.Range('E11').Value = '  -3.53%  '
$ws.Range('E12').Value = '  -3.08%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.900'
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('E14').Value = '  -6.46%  '
$ws.Range('D15').Value = '2.360.09'
$ws.Range('E15').Value = '  +0.38%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.36'
$ws.Range('E16').Value = '  -6.54%  '
$ws.Range('D17').Value = '2.069.80'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('D18').Value = '36.404.78'
$ws.Range('E18').Value = '  -1.80%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '16.67'
$ws.Range('E19').Value = '  -9.74%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '72.04'
$ws.Range('E20').Value = '  -3.41%  '
$ws.Range('D21').Value = '0.0₃0856'
$ws.Range('E21').Value = '  -5.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '238.04'
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.23'
$ws.Range('E23').Value = '  -4.47%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('E25').Value = '  -4.93%  '
$ws.Range('E26').Value = '  -1.88%  '
$ws.Range('E27').Value = '  -1.97%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '162.82'
$ws.Range('E28').Value = '  -4.44%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '20.15'
$ws.Range('E29').Value = '  +0.62%  '
$ws.Range('E30').Value = '  -2.05%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.20'
$ws.Range('E31').Value = '  +5.78%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.02'
$ws.Range('E32').Value = '  -6.76%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.46'
$ws.Range('E33').Value = '  -7.53%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0592'
$ws.Range('E34').Value = '  -4.46%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('E37').Value = '  -4.60%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0818'
$ws.Range('E38').Value = '  -6.72%  '
$ws.Range('E39').Value = '  -6.95%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.86'
$ws.Range('E40').Value = '  -5.37%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0215'
$ws.Range('E41').Value = '  -3.83%  '
$ws.Range('E42').Value = '  -8.46%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.10'
$ws.Range('E43').Value = '  -4.30%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '93.63'
$ws.Range('E44').Value = '  -5.13%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0903'
$ws.Range('E45').Value = '  -9.55%  '
$ws.Range('D46').Value = '1.395.40'
$ws.Range('E46').Value = '  +7.45%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '15.73'
$ws.Range('E47').Value = '  -8.40%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.44'
$ws.Range('E48').Value = '  +8.68%  '
$ws.Range('E49').Value = '  -0.66%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.26'
$ws.Range('E50').Value = '  -5.26%  '
$ws.Range('D51').Value = '2.246.31'
$ws.Range('E51').Value = '  +0.39%  '
